$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Arkansas) ---
$ws.Range("H5").Value = 5853
$ws.Range("I5").Value = 5819

# --- Row 6 ---
$ws.Range("AJ6").Value = 10164
$ws.Range("AK6").Value = 10048
$ws.Range("AL6").Value = 10097
$ws.Range("AO6").Value = 10680
$ws.Range("AP6").Value = 10755

# --- Row 7 ---
$ws.Range("AD7").Value = 5845
$ws.Range("AF7").Value = 4723
$ws.Range("AG7").Value = 4804

# --- Row 13 ---
$ws.Range("O13").Value = 41811
$ws.Range("V13").Value = 43664
$ws.Range("W13").Value = 44208
$ws.Range("X13").Value = 44840
$ws.Range("AA13").Value = 4329
$ws.Range("AB13").Value = 4329
$ws.Range("AC13").Value = 4258
$ws.Range("AD13").Value = 4622
$ws.Range("AE13").Value = 6582
$ws.Range("AF13").Value = 9380
$ws.Range("AG13").Value = 11762

# --- Row 14 ---
$ws.Range("N14").Value = 50987
$ws.Range("O14").Value = 51297
$ws.Range("P14").Value = 52164
$ws.Range("Q14").Value = 52784
$ws.Range("R14").Value = 53413
$ws.Range("S14").Value = 53858
$ws.Range("T14:V14").ClearContents()

# --- Row 15 ---
$ws.Range("J15").Value = 28871
$ws.Range("O15:V15").ClearContents()
$ws.Range("AG15").Value = 2981
$ws.Range("AP15").Value = 22957

# --- Row 18 ---
$ws.Range("AQ18").Value = 491000000

# --- Final selection (matches author's last-saved cursor position) ---
$ws.Range("AH7").Select() | Out-Null
